$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string JSON body for the new service row ---
$json = @"
{
    "success": true,
    "message": "Consulta exitosa.",
    "result": [
        {
            "id": 1,
            "clasifMovimientos": "COMPRA",
            "descripcionMovimientos": "COMPRA DE MERCANCIA"
        },
        {
            "id": 2,
            "clasifMovimientos": "VENTA",
            "descripcionMovimientos": "VENTA DE MERCANCIA"
        }
    ]
}
"@

# --- Remove the leftover yellow highlight from row 9 (A9:F9) ---
# (the service had previously been marked with a yellow fill; the commit
# clears that highlight while keeping the vertical-top / wrap alignment)
$ws.Range("A9:F9").Style = "Normal"
$ws.Range("A9:F9").VerticalAlignment = -4160
$ws.Range("B9").WrapText = $true
$ws.Range("F9").WrapText = $true

# --- Add the new "Listado tipo Movimiento" service as row 10 ---
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Listado tipo Movimiento"
$ws.Range("C10").Value = "http://localhost:8089/b-salesforce/rest/tipoMovimiento"
$ws.Range("D10").Value = "GET"
$ws.Range("F10").Value = $json

$ws.Range("A10:F10").VerticalAlignment = -4160
$ws.Range("B10").WrapText = $true
$ws.Range("F10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 270

# --- Update the view: scroll window down to the new row, select C18 ---
$ws.Range("C18").Select()
